# Applies:
#  1. Bold the "ID" (col A) and "Name" (col B) entries for the zone rows
#     (rows 4-9) on the "zones" sheet.
#  2. Re-seat the students within the "seating_plan_maths" sheet (same set
#     of names, new seat assignments) and move the active selection to K1.
#  3. Re-seat the students within the "seating_plan_english" sheet (same
#     set of names, new seat assignments).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "zones" sheet - bold the ID / Name columns for the zone definitions
# ---------------------------------------------------------------------
$wsZones = $wb.Worksheets.Item("zones")
$wsZones.Range("A4:B9").Font.Bold = $true

# ---------------------------------------------------------------------
# 2. "seating_plan_maths" sheet - shuffle students into new seats
# ---------------------------------------------------------------------
$wsMaths = $wb.Worksheets.Item("seating_plan_maths")

$wsMaths.Range("A2").Value = 'Brooke Layton, maths'
$wsMaths.Range("B2").Value = 'Stanley Hirst, maths'
$wsMaths.Range("C2").Value = 'James Eilbeck, maths'
$wsMaths.Range("D2").Value = 'Violet Hudson, maths'
$wsMaths.Range("E2").Value = 'Ruby Haigh, maths'
$wsMaths.Range("F2").Value = 'Madison Taylor, maths'
$wsMaths.Range("G2").Value = 'Sophie Rayner, maths'
$wsMaths.Range("H2").Value = 'Matthew Homan, maths'
$wsMaths.Range("I2").Value = 'Spencer Rowe, maths'
$wsMaths.Range("J2").Value = 'Niko Morris, maths'
$wsMaths.Range("B3").Value = 'James Shilton, maths'
$wsMaths.Range("C3").Value = 'Esther Sido, maths'
$wsMaths.Range("D3").Value = 'James Calderon, maths'
$wsMaths.Range("E3").Value = 'William Hunt, maths'
$wsMaths.Range("F3").Value = 'Katrina Petersone, maths'
$wsMaths.Range("G3").Value = 'Alex Sentance, maths'
$wsMaths.Range("H3").Value = 'Aarron Kelly, maths'
$wsMaths.Range("I3").Value = 'Lexi Green, maths'
$wsMaths.Range("J3").Value = 'Benjamin Hillary, maths'
$wsMaths.Range("B4").Value = 'Samuel Dixon, maths'
$wsMaths.Range("D4").Value = 'Ava Lee, maths'
$wsMaths.Range("E4").Value = 'Nancy Enyoazu, maths'
$wsMaths.Range("F4").Value = 'Caitlin Boyd, maths'
$wsMaths.Range("G4").Value = 'Benedict Hobday, maths'
$wsMaths.Range("H4").Value = 'Thomas Barrett, maths'

# ---------------------------------------------------------------------
# 3. "seating_plan_english" sheet - shuffle students into new seats
# ---------------------------------------------------------------------
$wsEnglish = $wb.Worksheets.Item("seating_plan_english")

$wsEnglish.Range("A2").Value = 'Hugo Bird, english'
$wsEnglish.Range("B2").Value = 'Patryk Rudnicki, english'
$wsEnglish.Range("C2").Value = 'Ava Lee, english'
$wsEnglish.Range("D2").Value = 'Lucy Webster, english'
$wsEnglish.Range("E2").Value = 'Jayden Parsons, english'
$wsEnglish.Range("F2").Value = 'Lewis Dacre, english'
$wsEnglish.Range("G2").Value = 'Cheryl Kanyimo, english'
$wsEnglish.Range("H2").Value = 'James Calderon, english'
$wsEnglish.Range("I2").Value = 'Lexie Starkey, english'
$wsEnglish.Range("J2").Value = 'Callum Foster, english'
$wsEnglish.Range("B3").Value = 'Niamh Teale, english'
$wsEnglish.Range("C3").Value = 'Matthew Homan, english'
$wsEnglish.Range("D3").Value = 'Ethan Durham, english'
$wsEnglish.Range("E3").Value = 'Caitlin Boyd, english'
$wsEnglish.Range("F3").Value = 'Samuel Dixon, english'
$wsEnglish.Range("G3").Value = 'Aarron Kelly, english'
$wsEnglish.Range("H3").Value = 'James Eilbeck, english'
$wsEnglish.Range("I3").Value = 'Mariam Keita, english'
$wsEnglish.Range("J3").Value = 'Elliott Long, english'
$wsEnglish.Range("B4").Value = 'Bethany Greer, english'
$wsEnglish.Range("C4").Value = 'Jayden Nasa-Mereni, english'
$wsEnglish.Range("D4").Value = 'Eva Redican, english'
$wsEnglish.Range("E4").Value = 'Cassie Strachan, english'
$wsEnglish.Range("F4").Value = 'Alex Sentance, english'
$wsEnglish.Range("G4").Value = 'Jude Fitzsimons, english'
$wsEnglish.Range("H4").Value = 'Isabella Holmes, english'

# ---------------------------------------------------------------------
# Move the active selection on "seating_plan_maths" to K1, then restore
# "seating_plan_english" as the active (tab-selected) sheet, matching the
# original workbook state.
# ---------------------------------------------------------------------
$wsMaths.Range("K1").Select()
$wsEnglish.Activate()
